# [Login] - BackEnd: Config for using Session/Cookie
#
# 1. The existing "Session/Cookie Authentication" reference link is swapped
#    out for a new c-sharpcorner article (both the visible text and the
#    underlying hyperlink target are updated).
# 2. A second bullet (same list level as the first link) is added right
#    after it, pointing to a companion c-sharpcorner article.

$d = $word.ActiveDocument

$oldUrlFragment = "tranvantoanblog.wordpress.com"
$newUrl1 = "https://www.c-sharpcorner.com/article/how-to-use-session-in-asp-net-core/?fbclid=IwAR0PgCrKsjsbzDpcK3NRdOOIC_-k3eEeLmGfrAyco0RRzdOje8fE0Tui-RQ"
$newUrl2 = "https://www.c-sharpcorner.com/article/all-about-session-in-asp-net-core/?fbclid=IwAR0k-UaN0ZES1oztyCbZzCRbDY7JZkDYhDBxBlZxtMaLKnBGL8Q20xutlOE"

# --- Step 1: locate the existing hyperlink that references the old blog post ---
$targetHyperlink = $null
for ($i = 1; $i -le $d.Hyperlinks.Count; $i++) {
    $candidate = $d.Hyperlinks.Item($i)
    if ($candidate.Address -like "*$oldUrlFragment*") {
        $targetHyperlink = $candidate
        break
    }
}

# --- Step 2: update its address and displayed text to the new article URL ---
$targetHyperlink.Address = $newUrl1
$targetHyperlink.TextToDisplay = $newUrl1

# --- Step 3: find the paragraph that now holds the updated hyperlink so we
#     can insert a new sibling bullet right after it ---
$hyperlinkParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*$newUrl1*") {
        $hyperlinkParaIndex = $i
        break
    }
}

$followingPara = $d.Paragraphs.Item($hyperlinkParaIndex + 1)

# Insert a brand-new empty paragraph right before the paragraph that follows
# the hyperlink bullet -- i.e. directly after the hyperlink bullet itself.
$insertionPoint = $followingPara.Range
$insertionPoint.Collapse(1)
$insertionPoint.InsertParagraphBefore()

# The freshly inserted paragraph is now at the same index the hyperlink's
# following paragraph used to occupy.
$newPara = $d.Paragraphs.Item($hyperlinkParaIndex + 1)

# Match the indentation level of the first link (w:ilvl 1 / second list level).
$newPara.Range.ListFormat.ListIndent()

# --- Step 4: put the new URL's text in the paragraph, then convert that
#     text (without its trailing paragraph mark) into a hyperlink ---
$newParaRange = $newPara.Range
$newParaRange.Text = $newUrl2

$newPara = $d.Paragraphs.Item($hyperlinkParaIndex + 1)
$textRange = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)
$d.Hyperlinks.Add($textRange, $newUrl2) | Out-Null
